# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column before the
# existing "Late" column (column N) so the schedule gains an extra
# "Variable Instalments" style column. This shifts the old Late / heading /
# Outstanding columns one position to the right (N->O, O->P, P->Q) and
# widens the new column to match its neighbour.
#
# Also move the active sheet / selection back to "Repayment schedule"
# (away from "Edit Repayment Schedule"), matching the state the workbook
# was left in after the edit.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Keep the new column's width consistent with its left neighbour (column M)
# instead of the default width.
$neighbourWidth = $wsSchedule.Columns("M").ColumnWidth

# Insert a blank column before column N ("Late"); existing N/O/P data
# (Late, heading/Over Due, Outstanding) shift right to O/P/Q.
$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = $neighbourWidth

# Restore the selection on "Edit Repayment Schedule" before leaving it, so
# its stored selection matches the post-edit state even though it is no
# longer the active sheet.
$wsEditSchedule = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEditSchedule.Range("D13").Select()

# Make "Repayment schedule" the active sheet again, with its own selection.
$wsSchedule.Activate()
$wsSchedule.Range("R7").Select()
